$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2) | Out-Null
}

# Title heading + bold byline near the end (same text, replaced globally)
Replace-Text "Play Country Farming for Free - Review" "Play Country Farming Free: Vibrant Slot Game with High Volatility"

# "What we like" bullet list
Replace-Text "High winning potential with 10,000x the bet payout" "Engaging gameplay with 20 fixed paylines"
Replace-Text "Vibrant and engaging cartoon-style graphics" "High volatility and theoretical RTP of 96.07%"
Replace-Text "Special symbols and features, including Wild and Scatter symbols" "Chance to win up to 10,000 times the bet"
Replace-Text "Multiple gameplay options for faster playtime" "Cartoon-style graphics with a pleasant farm life theme"

# "What we don't like" bullet list
Replace-Text "Only 20 fixed paylines" "Only one scatter symbol available for triggering free spins"
Replace-Text "Free spins can only be triggered by landing specific Scatter symbols" "Free spins can also be purchased at a cost of 100 times the bet"

# Meta description (italic run)
Replace-Text "Explore the features of Country Farming slot game with our review. Play for free and enjoy high winning potential, cartoon-style graphics, and special symbols." "Read our review of Country Farming slot game, play for free, and experience vibrant gameplay."
